$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look like plain numbers need an explicit
# Text number format first, otherwise Excel auto-converts the literal into
# a real number (e.g. "0.990" -> 0.99) and the trailing zero is lost.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated coin data (prices + 1h volume deltas); rows 48/49 also
# swap their Coin/Link content (WhiteBITCoin <-> RenderToken).
$ws.Range("D2").Value = '60.349.99'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '2.621.29'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '520.71'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").Value = '150.77'
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -4.14%  '
$ws.Range("D9").Value = '6.39'
$ws.Range("E9").Value = '  -4.47%  '
$ws.Range("E10").Value = '  +1.90%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").Value = '3.079.08'
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").Value = '60.384.08'
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").Value = '2.614.60'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").Value = '346.21'
$ws.Range("E19").Value = '  -3.21%  '
$ws.Range("D20").Value = '10.45'
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("D22").Value = '0.994'
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = '61.00'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("D25").Value = '0.163'
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").Value = '0.990'
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D27").Value = '0.0₃0834'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("E28").Value = '  -3.13%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '6.08'
$ws.Range("E30").Value = '  +2.66%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = '19.06'
$ws.Range("E32").Value = '  -1.96%  '
$ws.Range("D33").Value = '149.88'
$ws.Range("E33").Value = '  -0.41%  '
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("D35").Value = '0.898'
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("E36").Value = '  -2.40%  '
$ws.Range("D37").Value = '0.882'
$ws.Range("E37").Value = '  +4.47%  '
$ws.Range("D38").Value = '36.58'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("E40").Value = '  -1.94%  '
$ws.Range("D41").Value = '290.68'
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").Value = '0.630'
$ws.Range("E42").Value = '  +1.32%  '
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").Value = '0.0551'
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("D46").Value = '19.67'
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '4.75'
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = '10.35'
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("D50").Value = '18.97'
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("D51").Value = '1.966.98'
$ws.Range("E51").Value = '  -1.23%  '
